$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'personal pocket inhaler',
    'pickleball apparel women',
    'pickleball clothes',
    'pickleball clothes for women',
    'pickleball clothing women',
    'pickleball game',
    'pickleball gear',
    'pickleball now',
    'pickleball shorts',
    'pickleball shorts for women',
    'pickleball shorts women',
    'pickleball training',
    'pickleball training equipment',
    'pickleball women',
    'pickleball x',
    'picnic pants',
    'pilate pants',
    'pilate pants for women',
    'pilates ball base',
    'pilates ball with base',
    'pilates body ring',
    'pilates capris',
    'pilates chart',
    'pilates equipment capris',
    'pilates for runners',
    'pilates leggings for women',
    'pilates power ring',
    'pilates power ring pro',
    'pilates ring black',
    'pilates ring for women',
    'pilates ring pro',
    'pilates rings for women',
    'pilates training ring',
    'pilates workout chart',
    'pink and black tights',
    'pink compression leggings women',
    'pink compression tights',
    'pink tights for football',
    'pitch generator',
    'place sport pants',
    'plan b griptape',
    'plane care package',
    'plus compression leggings women',
    'plus knee brace for women',
    'plus knee high',
    'plus knee highs',
    'plus leggings capri',
    'plus size active wear women',
    'plus size athletic capris',
    'plus size basketball pants',
    'plus size black compression leggings',
    'plus size body tights',
    'plus size capri active leggings',
    'plus size capri leggings 2x',
    'plus size capri leggings for women',
    'plus size capri yoga pants',
    'plus size capris leggings',
    'plus size compression legging',
    'plus size compression leggings',
    'plus size compression leggings women',
    'plus size compression leggins',
    'plus size compression pants',
    'plus size compression pants women',
    'plus size compression running leggings',
    'plus size compression shorts women',
    'plus size compression tights',
    'plus size compression workout pants',
    'plus size fitness leggings',
    'plus size high stockings',
    'plus size high waist capri',
    'plus size high waist spanx',
    'plus size high waist workout leggings',
    'plus size knee high',
    'plus size knee high stockings',
    'plus size knee support',
    'plus size legginfs',
    'plus size legging tops',
    'plus size marijuana clothing',
    'plus size nike tights women',
    'plus size opaque tights',
    'plus size performance pants',
    'plus size running leggings',
    'plus size running pants',
    'plus size running tights',
    'plus size spanx',
    'plus size stocking',
    'plus size stockings',
    'plus size support hose',
    'plus size tights 3x',
    'plus size tights for women',
    'plus size tights high waist',
    'plus size tights knee high',
    'plus size white capri leggings',
    'plus size womens capri',
    'plus size womens compression shorts',
    'plus size womens running tights',
    'plus size workout',
    'plus size workout leggings',
    'plus size workout tights',
    'plus size workout tights for women'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
